$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.688.09'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.67%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.157.61'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.76%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '613.70'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.39'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.32%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.159.09'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.81%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.45'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.73%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.471'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.68%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000258'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.56'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.678.10'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.08%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.635.71'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.161.29'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.85'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '479.72'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.59'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.719'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.95'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +3.30%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.96'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.74%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.67'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.53%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.08'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.10%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.65%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -5.66%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.68'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.51'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.23%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.12'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0₃0786'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +8.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.99'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '53.22'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.18'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.96%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '460.35'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.08%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.70%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.31'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.860.32'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.31'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.267'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.98%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +6.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '26.54'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.17%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '35.69'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +9.11%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.50%  '
